# Add new ENGLISH quiz questions (rows 2-6) to the "ENGLISH" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENGLISH")

# Row data: Id, Subject Id, Question, Answer, Wrong Answers
$rows = @(
    @(9,  2, "___ you a student ?",             "Are",    "[Is, Do, Have]"),
    @(10, 2, "Translate 'book' into Uzbek",     "kitob",  "[qalam, ruchka, daftar]"),
    @(11, 2, "Translate 'apple' into Uzbek",    "olma",   "[gilos, olcha, anor]"),
    @(12, 2, "Translate 'dog' into Uzbek",      "it",     "[mushik, tovuq, sichqon]"),
    @(13, 2, "Translate 'cat' into Uzbek",      "mushuk", "[it, sichqon, tovuq]")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Resize columns to fit the new, wider content (matches target bestFit widths,
# expressed in "characters" as Excel's ColumnWidth expects: stored width - 5/6).
$ws.Columns.Item(1).ColumnWidth = 2.4361979166666665   # -> stored width 3.26953125
$ws.Columns.Item(3).ColumnWidth = 25.022135416666668   # -> stored width 25.85546875
$ws.Columns.Item(4).ColumnWidth = 7.291666666666667    # -> stored width 8.125
$ws.Columns.Item(5).ColumnWidth = 22.205729166666668   # -> stored width 23.0390625
